$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'track tights'
    2 = 'wintergear compression men'
    3 = 'basketball training tights'
    4 = 'basketball protector'
    5 = 'winter leggings for men'
    6 = 'kids tights with knee pads'
    7 = 'ropa de monta?a hombre'
    8 = 'men workout leggings nike'
    9 = 'boys white knee pads basketball'
    10 = 'coyote brown pants with knee pads'
    11 = 'black basketball knee pads nike'
    12 = 'basketball knee pads womens'
    13 = 'cold weather workout pants for men'
    14 = 'mens under armour long underwear pants'
    15 = 'calf compression leggings men'
    16 = 'tights for men nike'
    17 = 'swimming pants for men'
    18 = 'copper compression pants for men'
    19 = 'tesla thermal pants'
    20 = 'hex knee sleeve'
    21 = 'thermal nike'
    22 = 'running base layer men'
    23 = 'long spandex men'
    24 = 'cold gear compression pants men'
    25 = 'soccer winter gear'
    26 = 'kids soccor gear'
    27 = 'mens pants with knee pads'
    28 = 'mens thermal baselayer'
    29 = 'knee pad pants men'
    30 = 'adidas compression pants men'
    31 = 'adult clothing protector'
    32 = 'adult football knee pads'
    33 = 'adult football pants with pads'
    34 = 'athletic knee pads'
    35 = 'athletic tights men'
    36 = 'baseball knee pad'
    37 = 'baseball knee pads'
    38 = 'baseball pants adult small'
    39 = 'baskerball tights'
    40 = 'basketball compression knee pads'
    41 = 'basketball compression pants youth with knee pads'
    42 = 'basketball for youth'
    43 = 'basketball hip pads'
    44 = 'basketball leggings men'
    45 = 'basketball pants for women'
    46 = 'basketball pants men'
    47 = 'basketball tights boys youth'
    48 = 'basketball tights for men mcdavid'
    49 = 'basketball tights with pads for boys'
    50 = 'best basketball knee pads'
    51 = 'big knee pads'
    52 = 'big man knee pads'
    53 = 'bjj leggings'
    54 = 'black knee pads'
    55 = 'black knee pads for volleyball'
    56 = 'black leggings design'
    57 = 'black mesh leggings capri'
    58 = 'black youth baseball pants'
    59 = 'boys athletic tights basketball'
    60 = 'boys basketball knee pads mcdavid'
    61 = 'boys compression leggings'
    62 = 'boys compression leggings youth'
    63 = 'boys compression tights'
    64 = 'break away basketball pants'
    65 = 'capri leggings medium'
    66 = 'capri mens'
    67 = 'capri tights for men'
    68 = 'capris leggings'
    69 = 'cold gear for football'
    70 = 'cold weather panta'
    71 = 'compression for knee'
    72 = 'compression gear'
    73 = 'compression men pants'
    74 = 'compression pants big and tall men'
    75 = 'compression pants padded knees basketball'
    76 = 'compression shorts 3 4 length men'
    77 = 'compression tights for men'
    78 = 'cycling pants'
    79 = 'dark purple basketball knee pads'
    80 = 'elbow knee pads youth'
    81 = 'excersize equipment for men'
    82 = 'football 3 4 tights'
    83 = 'football knee pads'
    84 = 'football leg pads'
    85 = 'football pants'
    86 = 'football pants adult black'
    87 = 'football pants youth'
    88 = 'g form knee pads youth'
    89 = 'gel knee pads'
    90 = 'gel knee pads for men'
    91 = 'girl knee pads'
    92 = 'girls basketball knee pads'
    93 = 'girls tights with knee pads'
    94 = 'girls volleyball knee pads'
    95 = 'girls volleyball knee pads youth'
    96 = 'girls youth volleyball knee pads'
    97 = 'gym pants for men'
    98 = 'happy knees'
    99 = 'hex foam knee pads'
    100 = 'hex knee pads for basketball'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
